$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.692.19'
$ws.Range('E2').Value = '  +2.82%  '
$ws.Range('D3').Value = '2.326.37'
$ws.Range('E3').Value = '  +1.02%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '540.39'
$ws.Range('E5').Value = '  +4.63%  '
$ws.Range('D6').Value = '134.79'
$ws.Range('E6').Value = '  +3.39%  '
$ws.Range('E7').Value = '  +0.40%  '
$ws.Range('D8').Value = '0.561'
$ws.Range('E8').Value = '  +6.13%  '
$ws.Range('D9').Value = '0.101'
$ws.Range('E9').Value = '  +2.11%  '
$ws.Range('D10').Value = '5.46'
$ws.Range('E10').Value = '  +4.68%  '
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('E12').Value = '  +6.73%  '
$ws.Range('D13').Value = '23.66'
$ws.Range('E13').Value = '  +2.03%  '
$ws.Range('D14').Value = '2.748.10'
$ws.Range('E14').Value = '  +1.28%  '
$ws.Range('D15').Value = '57.670.98'
$ws.Range('E15').Value = '  +2.91%  '
$ws.Range('E16').Value = '  +1.38%  '
$ws.Range('D17').Value = '2.332.83'
$ws.Range('E17').Value = '  +1.16%  '
$ws.Range('D18').Value = '10.61'
$ws.Range('E18').Value = '  +3.07%  '
$ws.Range('D19').Value = '333.54'
$ws.Range('E19').Value = '  +1.56%  '
$ws.Range('D20').Value = '4.24'
$ws.Range('E20').Value = '  +3.03%  '
$ws.Range('D21').Value = '6.66'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').Value = '5.59'
$ws.Range('E23').Value = '  +0.96%  '
$ws.Range('D24').Value = '62.09'
$ws.Range('E24').Value = '  +2.02%  '
$ws.Range('D25').Value = '0.168'
$ws.Range('E25').Value = '  +2.86%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +1.05%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').Value = '8.46'
$ws.Range('E27').Value = '  -1.08%  '
$ws.Range('D28').Value = '1.39'
$ws.Range('E28').Value = '  +5.74%  '
$ws.Range('E29').Value = '  +4.39%  '
$ws.Range('D30').Value = '170.47'
$ws.Range('E30').Value = '  +2.05%  '
$ws.Range('D31').Value = '0.0₃0730'
$ws.Range('E31').Value = '  +2.82%  '
$ws.Range('D32').Value = '6.11'
$ws.Range('E32').Value = '  +1.30%  '
$ws.Range('E33').Value = '  +17.11%  '
$ws.Range('D34').Value = '18.43'
$ws.Range('E34').Value = '  +1.44%  '
$ws.Range('E36').Value = '  +0.54%  '
$ws.Range('E37').Value = '  +7.81%  '
$ws.Range('E38').Value = '  +1.94%  '
$ws.Range('E39').Value = '  +4.05%  '
$ws.Range('D40').Value = '39.11'
$ws.Range('E40').Value = '  +1.70%  '
$ws.Range('D41').Value = '144.83'
$ws.Range('E41').Value = '  -1.68%  '
$ws.Range('E42').Value = '  +0.50%  '
$ws.Range('E43').Value = '  +2.28%  '
$ws.Range('D44').Value = '284.39'
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('D45').Value = '0.0935'
$ws.Range('E45').Value = '  +1.42%  '
$ws.Range('D46').Value = '19.03'
$ws.Range('E46').Value = '  +5.29%  '
$ws.Range('D47').Value = '0.0501'
$ws.Range('E47').Value = '  +1.48%  '
$ws.Range('D48').Value = '0.558'
$ws.Range('E48').Value = '  +0.83%  '
$ws.Range('E49').Value = '  +2.34%  '
$ws.Range('D50').Value = '0.0215'
$ws.Range('E50').Value = '  +1.54%  '
$ws.Range('D51').Value = '17.41'
$ws.Range('E51').Value = '  +2.08%  '
